# v1.1 Changed the owner status
$wb = $excel.ActiveWorkbook

# --- Append a new version history entry ---
$historySheet = $wb.Worksheets.Item("Version History")
$historySheet.Range("D2").Formula = "=DATE(2025,4,21)"

$historySheet.Range("A3").Value = "v1.1"
$historySheet.Range("B3").Value = "Mahmoud Abdelmageed"
$historySheet.Range("C3").Value = "Changed the owner status"
$historySheet.Range("D3").Formula = "=DATE(2025,4,21)"

# --- Update Owner Status from Open to Closed on the notification reviews sheet ---
$reviewSheet = $wb.Worksheets.Item("LH_TC_NOTIFICATION_REVIEWS")
$reviewSheet.Range("I2").Value = "Closed"
$reviewSheet.Range("I3").Value = "Closed"
$reviewSheet.Range("I4").Value = "Closed"
